$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: paragraph "               - finished" (currently paragraph 6)
#   - remove the "_GoBack" bookmark that currently sits between the
#     two runs
#   - merge the leading "- " into the first run (spaces run) and
#     drop it from the second run so it reads "finished"
# ------------------------------------------------------------------

$finishedPara = $d.Paragraphs.Item(6)
$paraStart = $finishedPara.Range.Start

# First run: the leading spaces -> spaces + "- " (the bookmark that
# currently separates the two runs keeps them from merging back into
# a single run)
$spacesRange = $d.Range($paraStart, $paraStart + 15)
$spacesRange.Text = "               - "

# Second run: "- finished" -> "finished"
$tailStart = $spacesRange.End
$tailRange = $d.Range($tailStart, $tailStart + 10)
$tailRange.Text = "finished"

# Drop the old _GoBack bookmark from this paragraph now that the
# text is settled; it will be recreated later at the new edit
# location.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# Part 2: three new bulleted paragraphs after "Different Controlling
# Mechanisms ... March 31,2015"
# ------------------------------------------------------------------

$mechanismsPara = $d.Paragraphs.Item(9)
$insertPoint = $mechanismsPara.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$difficultyPara = $d.Paragraphs.Item(10)
$difficultyPara.Range.InsertBefore("Changing the difficulty level")

$rng2 = $d.Paragraphs.Item(10).Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

$cloudsPara = $d.Paragraphs.Item(11)
$cloudsPara.Range.InsertBefore("Making the clouds smaller as you go higher")
$d.Paragraphs.Item(11).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item(11).Range.Font.Bold = 1

$rng3 = $d.Paragraphs.Item(11).Range
$rng3.Collapse(0)
$rng3.InsertParagraphAfter()

$starsPara = $d.Paragraphs.Item(12)
$starsPara.Range.InsertBefore("Changing the stars")
$d.Paragraphs.Item(12).Range.ListFormat.ListLevelNumber = 2
$d.Paragraphs.Item(12).Range.Font.Bold = 1

# Re-create the "_GoBack" bookmark spanning the last two new bullets,
# mirroring where Word left it after this edit.
$cloudsPara = $d.Paragraphs.Item(11)
$starsPara = $d.Paragraphs.Item(12)
$newGoBackRange = $d.Range($cloudsPara.Range.Start, $starsPara.Range.End)
$d.Bookmarks.Add("_GoBack", $newGoBackRange)
